# Applies the per-row Coin/Link/Price/Volume(1h) updates described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") values are plain-looking decimals (e.g. "1.00", "0.999").
# A bare $cell.Value = "1.00" is auto-coerced by Excel to the number 1, losing the
# trailing zeros / text formatting the source data relies on (it is stored as text in
# the workbook, not a number). Forcing the cell to the "Text" number format for the
# duration of the write keeps the literal string, then restoring the "Normal" style
# afterwards drops the explicit style index again so the cell matches the original
# (unstyled) text cell.
function Set-TextValue {
    param($Cell, $Text)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "68.509.28"
$ws.Range("E2").Value = "  +0.91%  "

# Row 3
Set-TextValue $ws.Range("D3") "2.549.30"
$ws.Range("E3").Value = "  +0.70%  "

# Row 4
Set-TextValue $ws.Range("D4") "1.00"
$ws.Range("E4").Value = "  +0.14%  "

# Row 5
Set-TextValue $ws.Range("D5") "594.05"
$ws.Range("E5").Value = "  +0.18%  "

# Row 6
Set-TextValue $ws.Range("D6") "176.32"
$ws.Range("E6").Value = "  -0.11%  "

# Row 7
$ws.Range("E7").Value = "  +0.10%  "

# Row 8
Set-TextValue $ws.Range("D8") "0.526"
$ws.Range("E8").Value = "  -0.74%  "

# Row 9
Set-TextValue $ws.Range("D9") "2.548.81"
$ws.Range("E9").Value = "  +0.72%  "

# Row 10
$ws.Range("E10").Value = "  -1.50%  "

# Row 11
$ws.Range("E11").Value = "  +1.69%  "

# Row 12
$ws.Range("E12").Value = "  +0.62%  "

# Row 13
$ws.Range("E13").Value = "  -2.47%  "

# Row 14
Set-TextValue $ws.Range("D14") "26.65"
$ws.Range("E14").Value = "  -0.76%  "

# Row 15
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue $ws.Range("D15") "2.956.97"
$ws.Range("E15").Value = "  -0.87%  "

# Row 16
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue $ws.Range("D16") "0.0000178"
$ws.Range("E16").Value = "  -0.24%  "

# Row 17
Set-TextValue $ws.Range("D17") "68.522.22"
$ws.Range("E17").Value = "  +1.42%  "

# Row 18
Set-TextValue $ws.Range("D18") "2.22"
$ws.Range("E18").Value = "  +123.29%  "

# Row 19
Set-TextValue $ws.Range("D19") "2.542.55"
$ws.Range("E19").Value = "  +0.50%  "

# Row 20
Set-TextValue $ws.Range("D20") "11.94"
$ws.Range("E20").Value = "  +4.01%  "

# Row 21
$ws.Range("E21").Value = "  -0.60%  "

# Row 22
Set-TextValue $ws.Range("D22") "371.19"
$ws.Range("E22").Value = "  +3.15%  "

# Row 23
$ws.Range("E23").Value = "  -0.38%  "

# Row 24
Set-TextValue $ws.Range("D24") "4.59"
$ws.Range("E24").Value = "  -1.56%  "

# Row 25
Set-TextValue $ws.Range("D25") "71.94"
$ws.Range("E25").Value = "  +1.69%  "

# Row 26
$ws.Range("E26").Value = "  -0.02%  "

# Row 27
$ws.Range("E27").Value = "  -3.75%  "

# Row 28
$ws.Range("E28").Value = "  -3.04%  "

# Row 29
Set-TextValue $ws.Range("D29") "2.671.17"
$ws.Range("E29").Value = "  +0.51%  "

# Row 30
Set-TextValue $ws.Range("D30") "0.0₃0975"
$ws.Range("E30").Value = "  -1.40%  "

# Row 31
Set-TextValue $ws.Range("D31") "537.64"
$ws.Range("E31").Value = "  -2.77%  "

# Row 32
$ws.Range("E32").Value = "  +0.34%  "

# Row 33
Set-TextValue $ws.Range("D33") "1.32"
$ws.Range("E33").Value = "  -2.22%  "

# Row 34
$ws.Range("E34").Value = "  +0.69%  "

# Row 35
$ws.Range("E35").Value = "  -0.67%  "

# Row 36
Set-TextValue $ws.Range("D36") "0.999"
$ws.Range("E36").Value = "  -0.04%  "

# Row 37
Set-TextValue $ws.Range("D37") "160.43"
$ws.Range("E37").Value = "  +2.08%  "

# Row 38
$ws.Range("E38").Value = "  -2.09%  "

# Row 39
Set-TextValue $ws.Range("D39") "19.30"
$ws.Range("E39").Value = "  +2.76%  "

# Row 40
Set-TextValue $ws.Range("D40") "18.65"
$ws.Range("E40").Value = "  +0.23%  "

# Row 41
Set-TextValue $ws.Range("D41") "5.18"
$ws.Range("E41").Value = "  +0.07%  "

# Row 42
$ws.Range("E42").Value = "  -1.06%  "

# Row 43
Set-TextValue $ws.Range("D43") "0.351"
$ws.Range("E43").Value = "  -1.21%  "

# Row 44
$ws.Range("E44").Value = "  +0.17%  "

# Row 45
Set-TextValue $ws.Range("D45") "1.00"
$ws.Range("E45").Value = "  +0.39%  "

# Row 46
Set-TextValue $ws.Range("D46") "39.43"
$ws.Range("E46").Value = "  -1.31%  "

# Row 47
Set-TextValue $ws.Range("D47") "149.04"
$ws.Range("E47").Value = "  +0.54%  "

# Row 48
$ws.Range("E48").Value = "  +0.39%  "

# Row 49
$ws.Range("E49").Value = "  +0.73%  "

# Row 50
$ws.Range("E50").Value = "  -0.81%  "

# Row 51
Set-TextValue $ws.Range("D51") "1.72"
$ws.Range("E51").Value = "  +1.88%  "
